$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Apply the "filled" header/data style (currently used by Usuarios table, rows 10-11)
#        to the Inventario table (rows 7-8) and Detalle Ventas table (rows 13-14), which
#        previously used the non-filled variant.
$ws.Range("A10:F11").Copy() | Out-Null
$ws.Range("A7:F8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A10:F11").Copy() | Out-Null
$ws.Range("A13:F14").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$ws.Application.CutCopyMode = 0

# ... and to the "Ventas" table (rows 16-17), which only spans 4 columns (A:D).
$ws.Range("A10:C10").Copy() | Out-Null
$ws.Range("A16:C16").PasteSpecial(-4122) | Out-Null # xlPasteFormats (left + middle positions)
$ws.Range("F10").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null     # xlPasteFormats (right-most position)
$ws.Range("A11:D11").Copy() | Out-Null
$ws.Range("A17:D17").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$ws.Application.CutCopyMode = 0

# --- 2. Extend the "Inventario" title merge to include column F and re-merge.
$ws.Range("A7:E7").UnMerge() | Out-Null
$ws.Range("A7:F7").Merge() | Out-Null

# --- 3. Content fixes in the Inventario table.
$ws.Range("E8").Value = "fecha_inventario"
$ws.Range("F8").Value = "id_producto"

# --- 4. Fix a naming typo in Detalle Ventas (id_productos -> id_producto).
$ws.Range("F14").Value = "id_producto"

# --- 5. G13 / G14 (the lone extra column of the Detalle Ventas table) also switch from the
#        non-filled "box" style to the filled one, same as the rest of the table.
$ws.Range("A8").Copy() | Out-Null
$ws.Range("G13").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$ws.Range("G14").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$ws.Application.CutCopyMode = 0

# --- 6. View state tweaks (zoom + selection) to mirror the saved workbook view.
$ws.Application.ActiveWindow.Zoom = 65
$ws.Range("G22").Select() | Out-Null
